$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 14 rows before old row 47 so the "Spiel mit zwei Menschen"
# block shifts from rows 47-56 down to rows 61-70.
$ws.Rows("47:60").Insert()

# Rewrite the test-section header + body for rows 43-59
# ("Test: Doppelbesetzung" -> "Test: Spielfeldbesetzung - Randfaelle", expanded).
$ws.Range("A43").Value = "Test: Spielfeldbesetzung - Randfälle"
$ws.Range("A44").Value = "Der Nutzer wählt für Spieler 1 `"Mensch`", für Spieler 2 `"KI 1`" und klickt auf das Feld `"Spiel starten`"."
$ws.Range("B44").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an."
$ws.Range("A45").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben links."
$ws.Range("B45").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A46").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben rechts."
$ws.Range("B46").Value = "Es passiert nichts, da der Spieler `"Mensch`" nicht an der Reihe ist."
$ws.Range("A47").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B47").Value = "Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen `"O`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A48").Value = "Der Nutzer klickt auf dem Spielfeld auf das gerade von der KI ausgewählte Feld, welches mit einem `"O`" versehen wurde."
$ws.Range("B48").Value = "Es passiert nichts, da das Feld bereits besetzt ist."
$ws.Range("A49").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben links."
$ws.Range("B49").Value = "Es passiert nichts, da das Feld bereits besetzt ist."
$ws.Range("A50").Value = "Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf das `"NextMove`"-Symbol zu klicken, bis das Spiel vorbei ist und Spieler 1 gewonnen hat."
$ws.Range("B50").Value = "Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt. (Nicht implementiert: Es wird in die Belohnungsansicht gewechselt. Auf der rechten Seite wird als Graph der gesamte gewichtete Verlauf mit pro Zustand allen möglichen Äquivalenzklassenvertretern der Folgezuständen angezeigt.)"
$ws.Range("A51").Value = "Der Nutzer klickt auf das Feld `"Belohnung ausführen`"."
$ws.Range("B51").Value = "Es wird in den Startansicht gewechselt, in welchem die Spieler ausgewählt werden."
$ws.Range("A52").Value = "Der Nutzer wählt für Spieler 1 `"Mensch`", für Spieler 2 `"Mensch`" und klickt auf das Feld `"Spiel starten`"."
$ws.Range("B52").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an."
$ws.Range("A53").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben links."
$ws.Range("B53").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A54").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben links."
$ws.Range("B54").Value = "Es passiert nichts, da das Feld bereits besetzt ist."
$ws.Range("A55").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld unten links."
$ws.Range("B55").Value = "Das angeklickte Feld wird mit dem Zeichen `"O`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A56").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben in der Mitte."
$ws.Range("B56").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A57").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld unten in der Mitte."
$ws.Range("B57").Value = "Das angeklickte Feld wird mit dem Zeichen `"O`" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt."
$ws.Range("A58").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld oben rechts."
$ws.Range("B58").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird der aktualisierte Verlauf ohne Äquivalenzklassenvertreter der Folgezustände angezeigt. Das Spielergebnis wird angezeigt: `"Spieler 1 gewinnt!`""
$ws.Range("A59").Value = "Der Nutzer klickt auf dem Spielfeld auf das Feld unten rechts."
$ws.Range("B59").Value = "Es passiert nichts, da das Spiel bereits beendet ist."
$ws.Range("A61").Value = "Test: Spiel mit zwei Menschen"
$ws.Range("A62").Value = "Der Nutzer klickt auf den Button `"Neustart`"."
$ws.Range("B62").Value = "Die Spielerauswahl wird angezeigt."
$ws.Range("A63").Value = "Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 2."
$ws.Range("B63").Value = "Es wird im DropDown-Menü eine List aller möglichen Spieler angezeigt."
$ws.Range("A64").Value = "Der Nutzer wählt die Option `"Mensch`" im DropDown-Menü aus."
$ws.Range("B64").Value = "Als Spieler 2 wird ein Mensch festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde."
$ws.Range("A65").Value = "Der Nutzer klickt auf `"Spiel starten`""
$ws.Range("B65").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Folgezustände an."
$ws.Range("A66").Value = "Der Nutzer klickt auf das `"Play`"-Symbol."
$ws.Range("B66").Value = "Das `"Play`"-Symbol wird durch ein `"Pause`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."
$ws.Range("A67").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B67").Value = "Das `"Pause`"-Symbol wird durch ein `"Play`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."
$ws.Range("A68").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B68").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."
$ws.Range("A69").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B69").Value = "Das angeklickte Feld wird mit dem Zeichen `"O`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."
$ws.Range("A70").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B70").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

# A43 keeps the bold "section header" style (already bold from the original
# "Test: Doppelbesetzung" header); A61 keeps it too from the shifted original A47.
$ws.Range("A43").Font.Bold = $true
$ws.Range("A61").Font.Bold = $true

# Match the author re-selecting the last populated cell after the edit.
$ws.Range("A69").Select()
